$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 - Fi - 16.0 (1657)
$ws.Range("C3").Value = 677
$ws.Range("D3").Value = 81.5

# Row 4 - Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2
$ws.Range("D4").Value = 86.8

# Row 6 - Totals
$ws.Range("C6").Value = 1160

# Row 14 - Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("B14").Value = 449371

# Row 18 - Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("B18").Value = 77999
